$d = $word.ActiveDocument

# Shared envelope for the few structural tweaks that cannot be expressed as
# plain text edits (toggling the empty <w:lastRenderedPageBreak/> marker).
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXml($paragraphIndex, $innerParaXml) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $r = $p.Range
    $xml = $pkgOpen + $innerParaXml + $pkgClose
    $r.InsertXML($xml)
}

# --- 1. Drop the "kernel half-width must be scaled by pf" paragraph entirely ---
# It immediately follows the "w(x,y,z) = w(x) · w(y) · w(z)" paragraph.
$target = $d.Paragraphs.Item(12)
if ($target.Range.Text -notmatch "kernel half-width must be scaled by pf") {
    throw "Paragraph 12 text mismatch: $($target.Range.Text)"
}
$target.Range.Delete()

# NOTE: paragraph indices below are post-deletion (the removal above shifts
# every later paragraph's Paragraphs.Item index down by one).

# --- 2. Remove the stray <w:lastRenderedPageBreak/> before "Both operations..." ---
$p22 = $d.Paragraphs.Item(22)
if ($p22.Range.Text -notmatch "^Both operations evaluate weights") {
    throw "Paragraph 22 text mismatch: $($p22.Range.Text)"
}
Set-ParagraphXml 22 '<w:p w14:paraId="3AFC661F" w14:textId="77777777" w:rsidR="00204049" w:rsidRDefault="00BA6D34"><w:r><w:t>Both operations evaluate weights in padded logical units and operate on padded data samples. The target grid (polar or Cartesian) remains native.</w:t></w:r></w:p>'

# --- 3. Add <w:lastRenderedPageBreak/> to the "3. Polarization of 2D Particle Images" run ---
$p23 = $d.Paragraphs.Item(23)
if ($p23.Range.Text -notmatch "^3\. Polarization of 2D Particle Images") {
    throw "Paragraph 23 text mismatch: $($p23.Range.Text)"
}
Set-ParagraphXml 23 '<w:p w14:paraId="7FE698DC" w14:textId="77777777" w:rsidR="00204049" w:rsidRDefault="00BA6D34"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>3. Polarization of 2D Particle Images</w:t></w:r></w:p>'

# --- 4. Merge the split runs "...padded FFT lattice" + "." into a single run ---
$p27 = $d.Paragraphs.Item(27)
if ($p27.Range.Text -notmatch "^Interpolation is performed directly on the padded FFT lattice\.") {
    throw "Paragraph 27 text mismatch: $($p27.Range.Text)"
}
$p27.Range.Find.Execute("Interpolation is performed directly on the padded FFT lattice.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Interpolation is performed directly on the padded FFT lattice.", 2)

# --- 5. Merge the split runs "...Fourier volume" + "." into a single run ---
$p33 = $d.Paragraphs.Item(33)
if ($p33.Range.Text -notmatch "^Interpolation weights are evaluated in padded units at loc_pd\. Samples are gathered from the padded expanded Fourier volume\.") {
    throw "Paragraph 33 text mismatch: $($p33.Range.Text)"
}
$p33.Range.Find.Execute("Interpolation weights are evaluated in padded units at loc_pd. Samples are gathered from the padded expanded Fourier volume.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Interpolation weights are evaluated in padded units at loc_pd. Samples are gathered from the padded expanded Fourier volume.", 2)

Write-Output "Edits applied."
